# g 2019.05.22 No.03 modify
#
# 1. Paragraph 1 ("..."): drop the paragraph-mark run formatting
#    (<w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>).
# 2. Paragraph 2 ("..."): merge the two runs that used to be split by
#    the "_GoBack" bookmark into a single run, and give the paragraph
#    mark the eastAsia-hint run formatting that paragraph 1 used to have.
# 3. Append a brand-new paragraph ("和光同尘 gao") after paragraph 2,
#    keeping " gao" as its own run, and move the "_GoBack" bookmark to
#    sit at the end of this new paragraph.

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Paragraph 1: remove the eastAsia-hint paragraph-mark formatting ---
$para1Body = '<w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>谨以此书献给我的挚爱及亲人</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>!</w:t></w:r></w:p>'
$d.Paragraphs(1).Range.InsertXML($pkgHeader + $para1Body + $pkgFooter) | Out-Null

# --- Paragraph 2: single merged run + eastAsia-hint paragraph mark ---
$para2Body = '<w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>感谢生命中所有的遇见</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>!</w:t></w:r></w:p>'
$d.Paragraphs(2).Range.InsertXML($pkgHeader + $para2Body + $pkgFooter) | Out-Null

# --- New paragraph 3: "和光同尘" + " gao", with _GoBack at the end ---
$para3Body = '<w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>和光同尘</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> gao</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$endOfPara2 = $d.Paragraphs(2).Range.End
$d.Range($endOfPara2, $endOfPara2).InsertXML($pkgHeader + $para3Body + $pkgFooter) | Out-Null

Write-Output "done"
